# Generate Report for handback
#
# For each language sheet (zh-cn, de-de) mark the file as handed back:
#   - Status (col B)                 : "Ready for handoff" -> "Handed back: in sync with en-US"
#   - Latest Target File (col E)     : populate with the source file + hyperlink
#   - Latest Handback File (col F)   : populate with the handoff package + hyperlink
#   - Latest Handback DateTime (G)   : stamp with the handback time

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# ----- zh-cn -----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = $status

$ws.Hyperlinks.Add(
    $ws.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/ff4461a5faa889454f7a8977ee34728bc6d3dd2a/e2e/80044668-9a67-4ea8-bddf-41bd66cd9ed6.md",
    "",
    "",
    "80044668-9a67-4ea8-bddf-41bd66cd9ed6.md"
)

$ws.Hyperlinks.Add(
    $ws.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c770668245c79c278dbcb9b741046dc5e3e3337/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.zh-cn.xlf",
    "",
    "",
    "80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.zh-cn.xlf"
)

$ws.Range("G2").Value = "2016-01-14 02:22:51"

# ----- de-de -----
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = $status

$ws.Hyperlinks.Add(
    $ws.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/ff4461a5faa889454f7a8977ee34728bc6d3dd2a/e2e/80044668-9a67-4ea8-bddf-41bd66cd9ed6.md",
    "",
    "",
    "80044668-9a67-4ea8-bddf-41bd66cd9ed6.md"
)

$ws.Hyperlinks.Add(
    $ws.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7bbb1b353e91750f461af82b44c1d6a6fc92a581/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.de-de.xlf",
    "",
    "",
    "80044668-9a67-4ea8-bddf-41bd66cd9ed6.9b502d727bc474b6af29df6b8fbe69b87bf44395.de-de.xlf"
)

$ws.Range("G2").Value = "2016-01-14 02:23:14"

Write-Output "handback report generated"
